$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append the two new channel rows (rows 32-33) ---
# Row 32: "SUPER SOFIA IA" in both Channel (A) and CANAL (B) columns
$ws.Range("A32").Value = "SUPER SOFIA IA"
$ws.Range("B32").Value = "SUPER SOFIA IA"

# Row 33: "SOGNARE AGENTE IA" maps to "SUPER SOFIA IA"
$ws.Range("A33").Value = "SOGNARE AGENTE IA"
$ws.Range("B33").Value = "SUPER SOFIA IA"

# --- Match the recorded view state (best effort; some window-chrome
# metadata such as topLeftCell/zoomScaleNormal is not persisted by the
# host outside of freeze-pane scenarios, so this re-asserts intent
# without relying on it) ---
[void]$ws.Activate()
$excel.ActiveWindow.Zoom = 85
[void]$ws.Range("A16").Select()
$excel.ActiveWindow.ScrollRow = 16

# Final selection recorded in the workbook
[void]$ws.Range("F26").Select()
